{"js": "// Fix the typo \"form\" -> \"from\" in \"What can be learned form the data in\n// aggregate?\" (now reads \"... learned from the data in aggregate?\"), and\n// normalize a stray split-run (\"Have \" / \"groups\" / \" complete question 3.\")\n// back into a single run, matching the re-proofed / re-saved document.\n// Also relocate Word's \"_GoBack\" (last-edit) bookmark from its old spot to\n// the word that was just corrected.\n\nconst body = context.document.body;\n\n// --- 1) Merge the \"Have groups complete question 3.\" run split -----------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst haveGroupsText = \"Have groups complete question 3.\";\nlet haveGroupsPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === haveGroupsText) {\n    haveGroupsPara = paragraphs.items[i];\n    break;\n  }\n}\nif (haveGroupsPara) {\n  // Re-writing the full paragraph range collapses the old (proofing-split)\n  // runs into a single run with the same text.\n  haveGroupsPara.getRange().insertText(haveGroupsText, Word.InsertLocation.replace);\n}\n\n// --- 2) Fix the \"form\" -> \"from\" typo -------------------------------------\nconst needle = \"What can be learned form the data in aggregate?\";\nconst fixed = \"What can be learned from the data in aggregate?\";\n\nlet typoPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n    typoPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (typoPara) {\n  // Replacing the whole paragraph's text merges every run it contained\n  // into one run (preserving the shared bold/Arial formatting), the same\n  // way Word collapses runs after an edit + re-proof.\n  typoPara.getRange().insertText(fixed, Word.InsertLocation.replace);\n  await context.sync();\n\n  // --- 3) Move the \"_GoBack\" bookmark to the corrected word --------------\n  // Word always tracks the location of its most recent edit with a\n  // collapsed \"_GoBack\" bookmark. Drop the old one and re-create it right\n  // after the word that was just fixed.\n  const scoped = typoPara.getRange();\n  const hits = scoped.search(\"from\", { matchCase: true, matchWholeWord: true });\n  hits.load(\"text\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    const editPoint = hits.items[0].getRange(Word.RangeLocation.end);\n    context.document.deleteBookmark(\"_GoBack\");\n    editPoint.insertBookmark(\"_GoBack\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fixing a typo, form --> from\n#\n# \"What can be learned form the data in aggregate?\" becomes \"... learned\n# from the data in aggregate?\". Also normalizes a stray proofing-split run\n# (\"Have \" / \"groups\" / \" complete question 3.\") back into a single run,\n# and relocates Word's \"_GoBack\" (last-edit) bookmark from its old spot\n# onto the word that was just corrected - matching what Word itself does\n# when you make an edit and re-save.\n\n$d = $word.ActiveDocument\n\n# --- 1) Merge the \"Have groups complete question 3.\" run split -----------\n$haveGroupsText = \"Have groups complete question 3.\"\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Replacement.ClearFormatting()\n$r1.Find.Text = $haveGroupsText\n$r1.Find.Replacement.Text = $haveGroupsText\n$r1.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $haveGroupsText, 2) | Out-Null\n\n# --- 2) Fix the \"form\" -> \"from\" typo -------------------------------------\n$typoText = \"What can be learned form the data in aggregate?\"\n$fixedText = \"What can be learned from the data in aggregate?\"\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$r2.Find.Replacement.ClearFormatting()\n$r2.Find.Text = $typoText\n$r2.Find.Replacement.Text = $fixedText\n$r2.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $fixedText, 2) | Out-Null\n\n# --- 3) Move the \"_GoBack\" bookmark to the corrected word -----------------\n# Word always tracks the location of its most recent edit with a collapsed\n# \"_GoBack\" bookmark. Drop the old one and re-create it right after the\n# word that was just fixed.\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($fixedText)) {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -ne $null) {\n    $pr = $targetPara.Range\n    $localIdx = $pr.Text.IndexOf(\"from\")\n    if ($localIdx -ge 0) {\n        $editPoint = $pr.Start + $localIdx + 4\n        $newBookmarkRange = $d.Range($editPoint, $editPoint)\n\n        if ($d.Bookmarks.Exists(\"_GoBack\")) {\n            $d.Bookmarks.Item(\"_GoBack\").Delete()\n        }\n        $d.Bookmarks.Add(\"_GoBack\", $newBookmarkRange)\n    }\n}\n"}
